$d = $word.ActiveDocument

$replacements = @{
    "83×61=" = "46×61="
    "61×50=" = "28×26="
    "29×52=" = "84×20="
    "97×13=" = "47×25="
    "61×58=" = "63×74="
    "17×50=" = "48×14="
    "59×75=" = "98×60="
    "26×90=" = "18×86="
    "42×38=" = "39×28="
    "53×14=" = "63×83="
    "52×66=" = "15×38="
    "25×48=" = "82×18="
    "26×34=" = "86×84="
    "98×92=" = "45×38="
    "45×85=" = "27×70="
    "51×44=" = "23×93="
    "75×61=" = "42×37="
    "34×11=" = "26×45="
    "62×36=" = "50×39="
    "85×90=" = "60×92="
    "60×26=" = "34×43="
    "51×90=" = "72×20="
    "12×53=" = "53×44="
    "28×98=" = "60×97="
    "16×92=" = "21×51="
}

foreach ($old in $replacements.Keys) {
    $new = $replacements[$old]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
